$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new shared strings via new cell values
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B23").Value = "LOM3049 -  Termodinâmica de Máquinas  (Requisito)`n"
$ws.Range("C23").Value = "LOM3049 -  Termodinâmica de Máquinas  (Requisito)`n"

# Copy formatting from existing rows
$ws.Range("A17").Copy()
$ws.Range("A22").PasteSpecial(-4122)

$ws.Range("B21").Copy()
$ws.Range("B23").PasteSpecial(-4122)

$ws.Range("C21").Copy()
$ws.Range("C23").PasteSpecial(-4122)

$ws.Rows.Item(23).RowHeight = 30

Write-Output "done"
